# Tareas.xlsx edit script
# Adds four new task rows / updates row content according to the commit:
#   - "Que si cancelo imprimir al crear partida..." (already present, unchanged)
#   - "Agregar y quitar cotizaciones a Partidas ..." (replaces old row 4 content)
#   - "Agregar Inventarios a una Adquisicion ya creada, Fix" (new row)
#   - "Quitar Inventarios a una Adquisicion - Fix" (new row)
#   - "Documentos de Asignacion, rutas, idioma e imprimir" (new row, inserted
#     right after the existing "Documentos de rendicion..." row)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: make room for the 3 new rows that go right where the old
#     row 4 ("Documentos de rendicion...") used to live. That row's
#     content will end up shifted down to row 7.
$ws.Rows("4:6").Insert()

# New row 4: replaces the previous content, gets the green fill + wrap
# text style, a taller row and keeps B4 = 1.
$ws.Range("A4").Value2 = "Agregar y quitar cotizaciones a Partidas siempre tienen que quedar tres y el estado de SolicDetalle asociado no cambia"
$ws.Range("B4").Value2 = 1
$ws.Range("A4").Interior.Color = $ws.Range("A1").Interior.Color
$ws.Range("A4").WrapText = $true
$ws.Rows(4).RowHeight = 30

# New row 5: only column A is populated.
$ws.Range("A5").Value2 = "Agregar Inventarios a una Adquisición ya creada, Fix"
$ws.Range("A5").Interior.Color = $ws.Range("A1").Interior.Color
$ws.Range("A5").WrapText = $true

# New row 6: only column A is populated.
$ws.Range("A6").Value2 = "Quitar Inventarios a una Adquisicion - Fix"
$ws.Range("A6").Interior.Color = $ws.Range("A1").Interior.Color
$ws.Range("A6").WrapText = $true

# --- Step 2: insert one more row after the (now shifted) "Documentos de
#     rendicion..." row (row 7) for the new "Documentos de Asignacion..."
#     entry, which keeps the plain green-fill style used elsewhere (s=1).
$ws.Rows(8).Insert()
$ws.Range("A8").Value2 = "Documentos de Asignacion, rutas, idioma e imprimir"
$ws.Range("B8").Value2 = 2
$ws.Range("A8").Interior.Color = $ws.Range("A1").Interior.Color

# --- Step 3: update the selection to mirror the author's final cursor
#     position.
$ws.Range("A10").Select()

Write-Output "done"
